# Apply cryptocurrency price/volume updates to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "1.00" -> 1, dropping the fixed-decimal formatting the sheet relies on).
# Force them to Text format first, then restore the default style so no stray
# cell style is left behind, matching the un-styled inline-string cells in the source.
$textCells = 'D4', 'D5', 'D6', 'D10', 'D11', 'D16', 'D21', 'D22', 'D23', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D38', 'D39', 'D40', 'D41', 'D46', 'D47', 'D49', 'D50', 'D51'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '44.343.40'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.240.75'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '307.25'
$ws.Range('E5').Value = '  -2.97%  '
$ws.Range('D6').Value = '94.07'
$ws.Range('E6').Value = '  -6.00%  '
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('D10').Value = '34.70'
$ws.Range('E10').Value = '  -4.71%  '
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  -2.54%  '
$ws.Range('E12').Value = '  -4.00%  '
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = '2.337.48'
$ws.Range('E14').Value = '  +3.57%  '
$ws.Range('E15').Value = '  -2.73%  '
$ws.Range('D16').Value = '13.45'
$ws.Range('E16').Value = '  -4.34%  '
$ws.Range('D17').Value = '44.067.33'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '0.0₃0963'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('E20').Value = '  -9.42%  '
$ws.Range('D21').Value = '65.65'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '237.23'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = '2.94'
$ws.Range('E23').Value = '  -2.44%  '
$ws.Range('E24').Value = '  -2.44%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '40.02'
$ws.Range('E26').Value = '  +3.56%  '
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('D28').Value = '9.85'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('D29').Value = '5.91'
$ws.Range('E29').Value = '  -3.69%  '
$ws.Range('D30').Value = '20.02'
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').Value = '151.84'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').Value = '0.0797'
$ws.Range('E32').Value = '  -5.96%  '
$ws.Range('E33').Value = '  -2.90%  '
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  -12.64%  '
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('E37').Value = '  -9.54%  '
$ws.Range('D38').Value = '3.49'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').Value = '3.80'
$ws.Range('E39').Value = '  -4.83%  '
$ws.Range('D40').Value = '14.29'
$ws.Range('E40').Value = '  -7.94%  '
$ws.Range('D41').Value = '0.0298'
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = '1.701.28'
$ws.Range('E43').Value = '  -1.89%  '
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  -3.35%  '
$ws.Range('D46').Value = '99.54'
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').Value = '4.91'
$ws.Range('E47').Value = '  -5.87%  '
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '54.83'
$ws.Range('E49').Value = '  -4.34%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '8.07'
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('D51').Value = '67.24'
$ws.Range('E51').Value = '  -7.04%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
